$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Copy formatting from row 4 down to row 5 so the new row matches
# the style of the existing rows above it.
$ws.Range("A4:E4").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)

$ws.Range("A5").Value = 41699
$ws.Range("B5").Value = "03"
$ws.Range("C5").Value = "JMR"
$ws.Range("D5").Value = "Update from GPT naming to PIT"
$ws.Range("E5").Value = "Done"

$ws.Range("E5").Select()
